$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Decimal -> Hexadecimal/Binary conversion example: change decimal input from 100 to 10
$ws.Range("A3").Value = 10

# Hexadecimal calculator: change Operation from "-" to "+"
$ws.Range("A9").Value = "+"

# Change Operator 2 value from "d000" to "fc" (leading apostrophe keeps it
# stored as text with the original quote-prefix cell style, since this cell
# is number-formatted and would otherwise be reclassified)
$ws.Range("C9").Value = "'fc"

# Change Operator 1 value from "d800" to "fd"
$ws.Range("B9").Value = "fd"

$excel.CalculateFullRebuild()
